# Insert a new data row at row 9 (pushing existing rows 9..116 down to 10..117)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44552
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112045
$ws.Cells.Item(9, 7).Value = "Zapallo"
$ws.Cells.Item(9, 8).Value = "Camote"
$ws.Cells.Item(9, 9).Value = "1a nueva(o)"
$ws.Cells.Item(9, 10).Value = 600
$ws.Cells.Item(9, 11).Value = 700
$ws.Cells.Item(9, 12).Value = 750
$ws.Cells.Item(9, 13).Value = 725
$ws.Cells.Item(9, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(9, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 16).Value = 725
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
